$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 26, shifting rows 26:94 down to 27:95
$ws.Rows.Item(26).Insert()

# Populate the newly inserted row 26 with the new data record
$ws.Range("A26").Value = 9
$ws.Range("B26").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C26").Value = "Metropolitana"
$ws.Range("D26").Value = 44645
$ws.Range("D26").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E26").Value = 13
$ws.Range("F26").Value = "Fruta"
$ws.Range("G26").Value = 100101
$ws.Range("H26").Value = "Berries"
$ws.Range("I26").Value = 100101004
$ws.Range("J26").Value = "Frambuesa"
$ws.Range("K26").Value = "Sin especificar"
$ws.Range("L26").Value = "Primera"
$ws.Range("M26").Value = 450
$ws.Range("N26").Value = 8000
$ws.Range("O26").Value = 8000
$ws.Range("P26").Value = 8000
$ws.Range("Q26").Value = "$/bandeja 2 kilos"
$ws.Range("R26").Value = "Provincia de Linares"
$ws.Range("S26").Value = 4000
$ws.Range("T26").Value = 2
